$d = $word.ActiveDocument

# Locate the paragraph containing "Ver no Jupiter ..." (start of the block
# being removed) and the paragraph containing the "© 2020 ..." footer
# (end of the block being removed). The blank paragraph that separates
# "LOM3011: ..." from "Ver no Jupiter ..." is included in the deletion too,
# so we start the range at the end of the "LOM3011: ..." paragraph mark.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $p
    }
}

$start = $startPara.Range.Start
$end = $endPara.Range.End
$r = $d.Range($start, $end)
$r.Delete()
